# Update countries & provincias Spain
#
# Refreshes the COVID country data table on the "Pais" sheet:
#  - bumps the "last updated" timestamp (17:06 -> 18:23)
#  - updates the daily counters for several countries
#  - Mali/Cuba and Groenlandia/Islas Malvinas swap rank/position in the list

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1: last-updated banner
$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 18:23"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3644493
$ws.Range("C4").Value = 27666
$ws.Range("D4").Value = 1647032
$ws.Range("E4").Value = 1856984
$ws.Range("G4").Value = 333
$ws.Range("H4").Value = 140477

# Row 5: Brasil
$ws.Range("B5").Value = 1978236
$ws.Range("C5").Value = 7327
$ws.Range("E5").Value = 535764
$ws.Range("G5").Value = 174
$ws.Range("H5").Value = 75697

# Row 6: India
$ws.Range("B6").Value = 1001863
$ws.Range("C6").Value = 31694
$ws.Range("E6").Value = 356080
$ws.Range("G6").Value = 660
$ws.Range("H6").Value = 25589

# Row 16: Italia
$ws.Range("B16").Value = 243736
$ws.Range("C16").Value = 230
$ws.Range("D16").Value = 196246
$ws.Range("E16").Value = 12473
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = 35017

# Row 69: Chequia
$ws.Range("B69").Value = 13551
$ws.Range("C69").Value = 76
$ws.Range("E69").Value = 4689

# Row 101: Grecia
$ws.Range("B101").Value = 3939
$ws.Range("C101").Value = 29
$ws.Range("E101").Value = 2372

# Row 113: was Cuba, becomes Mali (new counters)
$ws.Range("A113").Value = "Mali"
$ws.Range("B113").Value = 2440
$ws.Range("C113").Value = 7
$ws.Range("D113").Value = 1777
$ws.Range("E113").Value = 542
$ws.Range("H113").Value = 121

# Row 114: was Mali, becomes Cuba (new counters)
$ws.Range("A114").Value = "Cuba"
$ws.Range("B114").Value = 2438
$ws.Range("D114").Value = 2277
$ws.Range("E114").Value = 74
$ws.Range("H114").Value = 87

# Row 119: Islandia
$ws.Range("B119").Value = 1914
$ws.Range("C119").Value = 3
$ws.Range("D119").Value = 1892
$ws.Range("E119").Value = 12

# Row 136: Jordania
$ws.Range("B136").Value = 1206
$ws.Range("C136").Value = 5
$ws.Range("D136").Value = 1019
$ws.Range("E136").Value = 177

# Row 159: Siria
$ws.Range("B159").Value = 477
$ws.Range("C159").Value = 19
$ws.Range("E159").Value = 315

# Row 209/210: Islas Malvinas <-> Groenlandia swap places (identical counters)
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"
